# Add a new Google-Forms testimonial response as row 10 of the
# "Form Responses 1" sheet, expand the "Form_Responses" table/autofilter
# and the _FilterDatabase defined name to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (row 9) down into row 10 so the
# new row inherits the exact same cell styles / number formats / row
# height as the other response rows.
$ws.Rows(9).Copy()
$ws.Rows(10).Insert()
$ws.Rows(10).RowHeight = 22.5

# Fill in the new testimonial's values.
$ws.Range("A10").Value = 46062.49713708334
$ws.Range("B10").Value = "lukumkulkarni@gmail.com"
$ws.Range("C10").Value = "Giving surface level answers that didn’t show the breadth of my knowledge. "
$ws.Range("D10").Value = "I wasn’t able to convey what was really going on in my projects and important my work was to the project. "
$ws.Range("E10").Value = "The most valuable piece of advice that you gave me was when you told me to try understand the core competency of the STAR question, rather than to loosely fitting my examples to the question. This way I can adapt my situation to specific topics that the interview is assessing, and I can bring out more detail about specifics that interview is inquiring about. "
$ws.Range("F10").Value = "I believe you are a no-nonsense mentor that brings what I need to work on straight to the point without sugarcoating. Sometimes other mentors will guide me in the wrong direction by not saying what I need to work on, because they are trying to not hurt feelings. But you are very direct and can pinpoint exactly what I need to work on, which streamlines the whole session. "
$ws.Range("G10").Value = "I was able to organize my thoughts by first understanding identifying the competency, then working backwards from there.  That even gave me better chances to pick a better story to satisfy the competency. "
$ws.Range("H10").Value = "My trajectory is upward since you helped me have more deliberate interview answers. I think this is the step in the right direction, but I need to do more work with you to solidify that trajectory. "
$ws.Range("I10").Value = "If you want someone to truly guide and asses how well you are doing without giving you sugarcoated feedback, then work with Mansour. He will give you his undivided attention as he deliberately finds ways you can improve your interviewing/career skills. "
$ws.Range("J10").Value = "Launch: getting the first job or an entry level job"
$ws.Range("K10").Value = "Yes"

# Expand the "Form_Responses" table (and its autofilter) to include the
# newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K10"))

# Update the _FilterDatabase defined name so it also covers row 10.
$n = $wb.Names.Item(1)
$n.RefersTo = "='Form Responses 1'!`$A`$1:`$K`$10"
